$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 40 (shifts old rows 40-49 down to 41-50)
$ws.Rows.Item(40).Insert()
$ws.Cells.Item(40, 1).Value = 108
$ws.Cells.Item(40, 2).Value = "Destocked - natural land"
$ws.Cells.Item(40, 4).Value = "Non-agricultural land-use"
$ws.Cells.Item(40, 5).Value = "Non-agricultural land-use"

# Insert two new rows at 47 (after the Biochar row, which is now row 46)
# This shifts old rows 47-50 (now) down to 49-52
$ws.Rows("47:48").Insert()

$ws.Cells.Item(47, 1).Value = 7
$ws.Cells.Item(47, 2).Value = "HIR - Beef"
$ws.Cells.Item(47, 5).Value = "Agricultural management"

$ws.Cells.Item(48, 1).Value = 8
$ws.Cells.Item(48, 2).Value = "HIR - Sheep"
$ws.Cells.Item(48, 5).Value = "Agricultural management"

# Update selection to match the target state
$ws.Range("B48").Select()
